$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$range = $ws.Range("A16:F20")
[void]$range.Select()
$range.ClearContents()
